# Team Meeting Attendance Form - add the two missed April meetings and
# drop the now-unused blank placeholder rows from the "Regularly Scheduled
# Meetings" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows 60/61 were blank placeholder rows (only B60/B61 had formatting, no
# data) directly below the last recorded meeting (row 59). Fill them in
# with the two additional regularly-scheduled "Team" meetings.
$ws.Range("A60").Value = "Team"
$ws.Range("B60").Value = "Monday, April 15,2019"
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 1

$ws.Range("A61").Value = "Team"
$ws.Range("B61").Value = "Wednesday, April 17, 2019"
$ws.Range("C61").Value = 1
$ws.Range("D61").Value = 1
$ws.Range("E61").Value = 1

# Rows 65-68 were left-over blank placeholder rows above the "Ad hoc
# Meetings" section; remove them entirely so the table closes up (all the
# totals/summary rows below shift up by 4 and their SUM ranges re-target
# automatically).
$ws.Rows("65:68").Delete()

# Restore the selection to the bottom of the (now shorter) sheet, matching
# where the last edit left the cursor.
$ws.Range("E83").Select() | Out-Null
